$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row2:
$ws.Range("I2").Value = "'79,976,616.95"
$ws.Range("J2").Value = "'79,145,809.76"
$ws.Range("K2").Value = "'80,376,068.20"
$ws.Range("L2").Value = "'79,868,341.35"
$ws.Range("M2").Value = "'78,634,868.31"
$ws.Range("N2").Value = "'80,254,117.59"
$ws.Range("O2").Value = "'79,393,109.68"
$ws.Range("P2").Value = "'78,815,218.77"
$ws.Range("Q2").Value = "'80,464,089.55"
$ws.Range("R2").Value = "'80,398,838.12"
$ws.Range("S2").Value = "'79,699,011.49"

# Row3:
$ws.Range("I3").Value = "'239,549,409.35"
$ws.Range("J3").Value = "'239,820,603.27"
$ws.Range("K3").Value = "'241,492,432.96"
$ws.Range("L3").Value = "'239,708,018.24"
$ws.Range("M3").Value = "'241,441,224.84"
$ws.Range("N3").Value = "'238,182,966.60"
$ws.Range("O3").Value = "'240,234,970.42"
$ws.Range("P3").Value = "'241,421,166.44"
$ws.Range("Q3").Value = "'239,436,983.89"
$ws.Range("R3").Value = "'239,448,383.43"
$ws.Range("S3").Value = "'238,536,048.93"

# Row4:
$ws.Range("I4").Value = "'242,004,004.60"
$ws.Range("J4").Value = "'239,898,473.56"
$ws.Range("K4").Value = "'239,008,833.36"
$ws.Range("L4").Value = "'237,795,169.88"
$ws.Range("M4").Value = "'238,906,165.35"
$ws.Range("N4").Value = "'240,697,047.84"
$ws.Range("O4").Value = "'239,381,295.54"
$ws.Range("P4").Value = "'238,493,110.07"
$ws.Range("Q4").Value = "'237,606,339.30"
$ws.Range("R4").Value = "'237,443,730.21"
$ws.Range("S4").Value = "'239,027,369.34"

# Row5:
$ws.Range("I5").Value = "'160,266,799.49"
$ws.Range("J5").Value = "'158,248,860.70"
$ws.Range("K5").Value = "'161,407,133.78"
$ws.Range("L5").Value = "'161,250,907.96"
$ws.Range("M5").Value = "'159,418,546.55"
$ws.Range("N5").Value = "'160,327,084.30"
$ws.Range("O5").Value = "'160,060,209.97"
$ws.Range("P5").Value = "'160,977,407.13"
$ws.Range("Q5").Value = "'161,936,277.26"
$ws.Range("R5").Value = "'160,055,872.03"
$ws.Range("S5").Value = "'160,104,398.94"
